# Commit: Thu, Jul 02, 2020 10:05:21 AM
#
# 1) Slide 16's table (3rd shape) switches from the custom "Table_0"
#    style ({126DB895-2345-447C-9A55-4874476DA5F2}) to the built-in
#    table style {1960CF9F-C41A-4B5A-9AFD-9C81C2546D8C}.
# 2) The deck's active theme (backing every slide/layout through the
#    single slide master) swaps its colour scheme from the "Integral"
#    palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -----------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}
$tableShape.Table.ApplyStyle("{1960CF9F-C41A-4B5A-9AFD-9C81C2546D8C}")

# --- 2. Theme colour scheme swap (Integral -> Office Theme) ------------
$clrScheme = $p.SlideMaster.Theme.ThemeColorScheme
$clrScheme.Item(1).RGB  = 0         # dk1      000000
$clrScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$clrScheme.Item(3).RGB  = 6968388   # dk2      44546A
$clrScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$clrScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$clrScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$clrScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$clrScheme.Item(8).RGB  = 49407     # accent4  FFC000
$clrScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$clrScheme.Item(10).RGB = 4697456   # accent6  70AD47
$clrScheme.Item(11).RGB = 12673797  # hlink    0563C1
$clrScheme.Item(12).RGB = 7491477   # folHlink 954F72
